$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "245.12"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "25.09"
$c.Style = "Normal"
$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "5.141"
$c.Style = "Normal"
$ws.Range("E4").Value = "3HuobiTokenHT"
$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.05635"
$c.Style = "Normal"
$ws.Range("E5").Value = "4CronosCRO"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "6.528"
$c.Style = "Normal"
$ws.Range("E6").Value = "5KuCoinTokenKCS"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "2.979"
$c.Style = "Normal"
$ws.Range("E7").Value = "6GateTokenGT"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.8118"
$c.Style = "Normal"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.8364"
$c.Style = "Normal"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.1336"
$c.Style = "Normal"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.06943"
$c.Style = "Normal"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.02843"
$c.Style = "Normal"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.09407"
$c.Style = "Normal"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001508"
$c.Style = "Normal"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0005986"
$c.Style = "Normal"
$ws.Range("E15").Value = "14OneONE"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.006116"
$c.Style = "Normal"
$ws.Range("E16").Value = "15TigerCashTCH"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.500"
$c.Style = "Normal"
$ws.Range("E17").Value = "16LEOLEO"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.03321"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.1291"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.766"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04687"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.1369"
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.001239"
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.004531"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.00009699"
$c.Style = "Normal"
$ws.Range("E27").Value = "26NitroExNTXBestin24h"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0001868"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.03622"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.006239"
$c.Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1051"
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.002716"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.008328"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00005280"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.2199"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.002286"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.00002099"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0001999"
$c.Style = "Normal"
